$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Fecha) - cyclic shift of date serials across rows 2-5
$ws.Range("D2").Value = 44250
$ws.Range("D3").Value = 44257
$ws.Range("D4").Value = 44252
$ws.Range("D5").Value = 44253

# Column M (Volumen)
$ws.Range("M2").Value = 200
$ws.Range("M3").Value = 100
$ws.Range("M4").Value = 120
$ws.Range("M5").Value = 160

# Column N (Precio minimo)
$ws.Range("N4").Value = 13000
$ws.Range("N5").Value = 14000

# Column O (Precio maximo)
$ws.Range("O4").Value = 14000
$ws.Range("O5").Value = 15000

# Column P (Precio promedio ponderado)
$ws.Range("P4").Value = 13500
$ws.Range("P5").Value = 14500

# Column S (Precio $/Kg)
$ws.Range("S4").Value = 750
$ws.Range("S5").Value = 806
